$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "索引" (Index) column F is no longer needed - remove it entirely.
# This shifts everything in columns G:S left by one (to F:R).
$ws.Columns.Item(6).Delete()

# Fix the example value for the "自动订阅" (Auto Subscribe) column - the
# sample row had an ambiguous "0/1" placeholder; replace it with a single
# concrete example value "0".
$ws.Range("R3").Value = "0"
